# =====================================================================
# Commit: "Add files via upload"
# Adds a new sheet "Plan3" (becomes the active tab) containing a 21x21
# municipality adjacency/co-occurrence matrix, and updates the saved
# cursor position / selection on all three sheets.
# =====================================================================

$wb = $excel.ActiveWorkbook

# --- 1. Update Plan1 (sheet1) view: scroll to A35, select A51 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Application.Goto($ws1.Range("A35"), $false)
$ws1.Range("A51").Select()

# --- 2. Update Plan2 (sheet2) view: scroll to A1 (no frozen/offset), select A4 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Application.Goto($ws2.Range("A1"), $false)
$ws2.Range("A4").Select()

# --- 3. Add Plan3 as a new sheet after Plan2; it becomes the active sheet/tab ---
$lastIndex = $wb.Worksheets.Count
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$ws3.Name = "Plan3"

# --- 4. Column widths (Plan3) ---
$ws3.Columns.Item(1).ColumnWidth = 29.0
$ws3.Columns.Item(2).ColumnWidth = 12.166666666666666
$ws3.Columns.Item(3).ColumnWidth = 11.666666666666666
$ws3.Columns.Item(4).ColumnWidth = 15.5
$ws3.Columns.Item(5).ColumnWidth = 10.166666666666666
$ws3.Columns.Item(6).ColumnWidth = 13.666666666666666
$ws3.Columns.Item(7).ColumnWidth = 14.0
$ws3.Columns.Item(8).ColumnWidth = 10.666666666666666
$ws3.Columns.Item(9).ColumnWidth = 28.333333333333332
$ws3.Columns.Item(10).ColumnWidth = 28.333333333333332
$ws3.Columns.Item(11).ColumnWidth = 27.5
$ws3.Columns.Item(12).ColumnWidth = 28.833333333333332
$ws3.Columns.Item(13).ColumnWidth = 11.833333333333334
$ws3.Columns.Item(15).ColumnWidth = 11.833333333333334
$ws3.Columns.Item(16).ColumnWidth = 14.333333333333334
$ws3.Columns.Item(17).ColumnWidth = 16.5
$ws3.Columns.Item(18).ColumnWidth = 16.0
$ws3.Columns.Item(19).ColumnWidth = 8.666666666666666
$ws3.Columns.Item(20).ColumnWidth = 15.166666666666666
$ws3.Columns.Item(21).ColumnWidth = 13.166666666666666
$ws3.Columns.Item(22).ColumnWidth = 11.666666666666666

# --- 5. Header row (B1:V1): municipality names ---
$ws3.Range("B1").Value = "Itabaiana"
$ws3.Range("C1").Value = "Frei Paulo"
$ws3.Range("D1").Value = "Moita Bonita"
$ws3.Range("E1").Value = "Lagarto"
$ws3.Range("F1").Value = "Laranjeiras"
$ws3.Range("G1").Value = "Macambira"
$ws3.Range("H1").Value = "Malhador"
$ws3.Range("I1").Value = "Nossa Senhora Aparecida"
$ws3.Range("J1").Value = "Nossa Senhora da Gloria"
$ws3.Range("K1").Value = "Nossa Senhora das Dores"
$ws3.Range("L1").Value = "Nossa Senhora do Socorro"
$ws3.Range("M1").Value = "Pedra Mole"
$ws3.Range("N1").Value = "Pinhão"
$ws3.Range("O1").Value = "Riachuelo"
$ws3.Range("P1").Value = "Ribeirópolis"
$ws3.Range("Q1").Value = "São Cristóvão"
$ws3.Range("R1").Value = "São Domingos"
$ws3.Range("S1").Value = "Carira"
$ws3.Range("T1").Value = "Areia Branca"
$ws3.Range("U1").Value = "Simão Dias"
$ws3.Range("V1").Value = "Aracaju"

# --- 6. Data grid (A2:V22): row label + 21x21 0/1 matrix ---
$ws3.Range("A2").Value = "Itabaiana"
$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = 1
$ws3.Range("D2").Value = 1
$ws3.Range("E2").Value = 0
$ws3.Range("F2").Value = 0
$ws3.Range("G2").Value = 1
$ws3.Range("H2").Value = 1
$ws3.Range("I2").Value = 0
$ws3.Range("J2").Value = 0
$ws3.Range("K2").Value = 0
$ws3.Range("L2").Value = 0
$ws3.Range("M2").Value = 0
$ws3.Range("N2").Value = 0
$ws3.Range("O2").Value = 0
$ws3.Range("P2").Value = 1
$ws3.Range("Q2").Value = 0
$ws3.Range("R2").Value = 0
$ws3.Range("S2").Value = 0
$ws3.Range("T2").Value = 1
$ws3.Range("U2").Value = 0
$ws3.Range("V2").Value = 0
$ws3.Range("A3").Value = "Frei Paulo"
$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = 0
$ws3.Range("D3").Value = 1
$ws3.Range("E3").Value = 0
$ws3.Range("F3").Value = 0
$ws3.Range("G3").Value = 1
$ws3.Range("H3").Value = 0
$ws3.Range("I3").Value = 1
$ws3.Range("J3").Value = 1
$ws3.Range("K3").Value = 0
$ws3.Range("L3").Value = 0
$ws3.Range("M3").Value = 1
$ws3.Range("N3").Value = 1
$ws3.Range("O3").Value = 0
$ws3.Range("P3").Value = 1
$ws3.Range("Q3").Value = 0
$ws3.Range("R3").Value = 0
$ws3.Range("S3").Value = 0
$ws3.Range("T3").Value = 0
$ws3.Range("U3").Value = 0
$ws3.Range("V3").Value = 0
$ws3.Range("A4").Value = "Moita Bonita"
$ws3.Range("B4").Value = 1
$ws3.Range("C4").Value = 0
$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 0
$ws3.Range("F4").Value = 0
$ws3.Range("G4").Value = 0
$ws3.Range("H4").Value = 1
$ws3.Range("I4").Value = 0
$ws3.Range("J4").Value = 0
$ws3.Range("K4").Value = 1
$ws3.Range("L4").Value = 0
$ws3.Range("M4").Value = 0
$ws3.Range("N4").Value = 0
$ws3.Range("O4").Value = 0
$ws3.Range("P4").Value = 1
$ws3.Range("Q4").Value = 0
$ws3.Range("R4").Value = 0
$ws3.Range("S4").Value = 0
$ws3.Range("T4").Value = 0
$ws3.Range("U4").Value = 0
$ws3.Range("V4").Value = 0
$ws3.Range("A5").Value = "Lagarto"
$ws3.Range("B5").Value = 0
$ws3.Range("C5").Value = 0
$ws3.Range("D5").Value = 0
$ws3.Range("E5").Value = 0
$ws3.Range("F5").Value = 0
$ws3.Range("G5").Value = 0
$ws3.Range("H5").Value = 0
$ws3.Range("I5").Value = 0
$ws3.Range("J5").Value = 0
$ws3.Range("K5").Value = 0
$ws3.Range("L5").Value = 0
$ws3.Range("M5").Value = 0
$ws3.Range("N5").Value = 1
$ws3.Range("O5").Value = 0
$ws3.Range("P5").Value = 0
$ws3.Range("Q5").Value = 0
$ws3.Range("R5").Value = 1
$ws3.Range("S5").Value = 0
$ws3.Range("T5").Value = 0
$ws3.Range("U5").Value = 1
$ws3.Range("V5").Value = 0
$ws3.Range("A6").Value = "Laranjeiras"
$ws3.Range("B6").Value = 0
$ws3.Range("C6").Value = 1
$ws3.Range("D6").Value = 0
$ws3.Range("E6").Value = 0
$ws3.Range("F6").Value = 0
$ws3.Range("G6").Value = 0
$ws3.Range("H6").Value = 0
$ws3.Range("I6").Value = 0
$ws3.Range("J6").Value = 0
$ws3.Range("K6").Value = 0
$ws3.Range("L6").Value = 0
$ws3.Range("M6").Value = 0
$ws3.Range("N6").Value = 0
$ws3.Range("O6").Value = 1
$ws3.Range("P6").Value = 0
$ws3.Range("Q6").Value = 1
$ws3.Range("R6").Value = 0
$ws3.Range("S6").Value = 0
$ws3.Range("T6").Value = 1
$ws3.Range("U6").Value = 0
$ws3.Range("V6").Value = 1
$ws3.Range("A7").Value = "Macambira"
$ws3.Range("B7").Value = 1
$ws3.Range("C7").Value = 1
$ws3.Range("D7").Value = 0
$ws3.Range("E7").Value = 0
$ws3.Range("F7").Value = 0
$ws3.Range("G7").Value = 0
$ws3.Range("H7").Value = 0
$ws3.Range("I7").Value = 0
$ws3.Range("J7").Value = 0
$ws3.Range("K7").Value = 0
$ws3.Range("L7").Value = 0
$ws3.Range("M7").Value = 1
$ws3.Range("N7").Value = 1
$ws3.Range("O7").Value = 0
$ws3.Range("P7").Value = 1
$ws3.Range("Q7").Value = 0
$ws3.Range("R7").Value = 1
$ws3.Range("S7").Value = 1
$ws3.Range("T7").Value = 0
$ws3.Range("U7").Value = 0
$ws3.Range("V7").Value = 0
$ws3.Range("A8").Value = "Malhador"
$ws3.Range("B8").Value = 1
$ws3.Range("C8").Value = 0
$ws3.Range("D8").Value = 1
$ws3.Range("E8").Value = 0
$ws3.Range("F8").Value = 0
$ws3.Range("G8").Value = 0
$ws3.Range("H8").Value = 0
$ws3.Range("I8").Value = 0
$ws3.Range("J8").Value = 0
$ws3.Range("K8").Value = 0
$ws3.Range("L8").Value = 0
$ws3.Range("M8").Value = 0
$ws3.Range("N8").Value = 0
$ws3.Range("O8").Value = 1
$ws3.Range("P8").Value = 0
$ws3.Range("Q8").Value = 0
$ws3.Range("R8").Value = 0
$ws3.Range("S8").Value = 0
$ws3.Range("T8").Value = 1
$ws3.Range("U8").Value = 0
$ws3.Range("V8").Value = 0
$ws3.Range("A9").Value = "Nossa Senhora Aparecida"
$ws3.Range("B9").Value = 0
$ws3.Range("C9").Value = 1
$ws3.Range("D9").Value = 0
$ws3.Range("E9").Value = 0
$ws3.Range("F9").Value = 0
$ws3.Range("G9").Value = 0
$ws3.Range("H9").Value = 0
$ws3.Range("I9").Value = 0
$ws3.Range("J9").Value = 1
$ws3.Range("K9").Value = 0
$ws3.Range("L9").Value = 0
$ws3.Range("M9").Value = 1
$ws3.Range("N9").Value = 1
$ws3.Range("O9").Value = 0
$ws3.Range("P9").Value = 1
$ws3.Range("Q9").Value = 0
$ws3.Range("R9").Value = 0
$ws3.Range("S9").Value = 1
$ws3.Range("T9").Value = 0
$ws3.Range("U9").Value = 0
$ws3.Range("V9").Value = 0
$ws3.Range("A10").Value = "Nossa Senhora da Gloria"
$ws3.Range("B10").Value = 0
$ws3.Range("C10").Value = 1
$ws3.Range("D10").Value = 0
$ws3.Range("E10").Value = 0
$ws3.Range("F10").Value = 0
$ws3.Range("G10").Value = 0
$ws3.Range("H10").Value = 0
$ws3.Range("I10").Value = 1
$ws3.Range("J10").Value = 0
$ws3.Range("K10").Value = 0
$ws3.Range("L10").Value = 0
$ws3.Range("M10").Value = 1
$ws3.Range("N10").Value = 1
$ws3.Range("O10").Value = 0
$ws3.Range("P10").Value = 1
$ws3.Range("Q10").Value = 0
$ws3.Range("R10").Value = 0
$ws3.Range("S10").Value = 1
$ws3.Range("T10").Value = 0
$ws3.Range("U10").Value = 0
$ws3.Range("V10").Value = 0
$ws3.Range("A11").Value = "Nossa Senhora das Dores"
$ws3.Range("B11").Value = 0
$ws3.Range("C11").Value = 0
$ws3.Range("D11").Value = 1
$ws3.Range("E11").Value = 0
$ws3.Range("F11").Value = 0
$ws3.Range("G11").Value = 0
$ws3.Range("H11").Value = 0
$ws3.Range("I11").Value = 0
$ws3.Range("J11").Value = 0
$ws3.Range("K11").Value = 0
$ws3.Range("L11").Value = 0
$ws3.Range("M11").Value = 0
$ws3.Range("N11").Value = 0
$ws3.Range("O11").Value = 0
$ws3.Range("P11").Value = 1
$ws3.Range("Q11").Value = 0
$ws3.Range("R11").Value = 0
$ws3.Range("S11").Value = 0
$ws3.Range("T11").Value = 0
$ws3.Range("U11").Value = 0
$ws3.Range("V11").Value = 0
$ws3.Range("A12").Value = "Nossa Senhora do Socorro"
$ws3.Range("B12").Value = 0
$ws3.Range("C12").Value = 0
$ws3.Range("D12").Value = 0
$ws3.Range("E12").Value = 0
$ws3.Range("F12").Value = 1
$ws3.Range("G12").Value = 0
$ws3.Range("H12").Value = 0
$ws3.Range("I12").Value = 0
$ws3.Range("J12").Value = 0
$ws3.Range("K12").Value = 0
$ws3.Range("L12").Value = 0
$ws3.Range("M12").Value = 0
$ws3.Range("N12").Value = 0
$ws3.Range("O12").Value = 0
$ws3.Range("P12").Value = 0
$ws3.Range("Q12").Value = 0
$ws3.Range("R12").Value = 0
$ws3.Range("S12").Value = 0
$ws3.Range("T12").Value = 0
$ws3.Range("U12").Value = 0
$ws3.Range("V12").Value = 1
$ws3.Range("A13").Value = "Pedra Mole"
$ws3.Range("B13").Value = 0
$ws3.Range("C13").Value = 1
$ws3.Range("D13").Value = 0
$ws3.Range("E13").Value = 0
$ws3.Range("F13").Value = 0
$ws3.Range("G13").Value = 1
$ws3.Range("H13").Value = 0
$ws3.Range("I13").Value = 1
$ws3.Range("J13").Value = 1
$ws3.Range("K13").Value = 0
$ws3.Range("L13").Value = 0
$ws3.Range("M13").Value = 0
$ws3.Range("N13").Value = 1
$ws3.Range("O13").Value = 0
$ws3.Range("P13").Value = 1
$ws3.Range("Q13").Value = 0
$ws3.Range("R13").Value = 0
$ws3.Range("S13").Value = 1
$ws3.Range("T13").Value = 0
$ws3.Range("U13").Value = 1
$ws3.Range("V13").Value = 0
$ws3.Range("A14").Value = "Pinhão"
$ws3.Range("B14").Value = 0
$ws3.Range("C14").Value = 1
$ws3.Range("D14").Value = 0
$ws3.Range("E14").Value = 1
$ws3.Range("F14").Value = 0
$ws3.Range("G14").Value = 1
$ws3.Range("H14").Value = 0
$ws3.Range("I14").Value = 1
$ws3.Range("J14").Value = 1
$ws3.Range("K14").Value = 0
$ws3.Range("L14").Value = 0
$ws3.Range("M14").Value = 1
$ws3.Range("N14").Value = 0
$ws3.Range("O14").Value = 0
$ws3.Range("P14").Value = 1
$ws3.Range("Q14").Value = 0
$ws3.Range("R14").Value = 0
$ws3.Range("S14").Value = 1
$ws3.Range("T14").Value = 0
$ws3.Range("U14").Value = 1
$ws3.Range("V14").Value = 0
$ws3.Range("A15").Value = "Riachuelo"
$ws3.Range("B15").Value = 0
$ws3.Range("C15").Value = 0
$ws3.Range("D15").Value = 0
$ws3.Range("E15").Value = 0
$ws3.Range("F15").Value = 1
$ws3.Range("G15").Value = 0
$ws3.Range("H15").Value = 1
$ws3.Range("I15").Value = 0
$ws3.Range("J15").Value = 0
$ws3.Range("K15").Value = 0
$ws3.Range("L15").Value = 0
$ws3.Range("M15").Value = 0
$ws3.Range("N15").Value = 0
$ws3.Range("O15").Value = 0
$ws3.Range("P15").Value = 0
$ws3.Range("Q15").Value = 0
$ws3.Range("R15").Value = 0
$ws3.Range("S15").Value = 0
$ws3.Range("T15").Value = 1
$ws3.Range("U15").Value = 0
$ws3.Range("V15").Value = 1
$ws3.Range("A16").Value = "Ribeirópolis"
$ws3.Range("B16").Value = 1
$ws3.Range("C16").Value = 1
$ws3.Range("D16").Value = 1
$ws3.Range("E16").Value = 0
$ws3.Range("F16").Value = 0
$ws3.Range("G16").Value = 1
$ws3.Range("H16").Value = 0
$ws3.Range("I16").Value = 1
$ws3.Range("J16").Value = 1
$ws3.Range("K16").Value = 1
$ws3.Range("L16").Value = 0
$ws3.Range("M16").Value = 1
$ws3.Range("N16").Value = 1
$ws3.Range("O16").Value = 0
$ws3.Range("P16").Value = 0
$ws3.Range("Q16").Value = 0
$ws3.Range("R16").Value = 0
$ws3.Range("S16").Value = 1
$ws3.Range("T16").Value = 0
$ws3.Range("U16").Value = 0
$ws3.Range("V16").Value = 0
$ws3.Range("A17").Value = "São Cristóvão"
$ws3.Range("B17").Value = 0
$ws3.Range("C17").Value = 0
$ws3.Range("D17").Value = 0
$ws3.Range("E17").Value = 0
$ws3.Range("F17").Value = 1
$ws3.Range("G17").Value = 0
$ws3.Range("H17").Value = 0
$ws3.Range("I17").Value = 0
$ws3.Range("J17").Value = 0
$ws3.Range("K17").Value = 1
$ws3.Range("L17").Value = 1
$ws3.Range("M17").Value = 0
$ws3.Range("N17").Value = 0
$ws3.Range("O17").Value = 1
$ws3.Range("P17").Value = 0
$ws3.Range("Q17").Value = 0
$ws3.Range("R17").Value = 0
$ws3.Range("S17").Value = 0
$ws3.Range("T17").Value = 1
$ws3.Range("U17").Value = 0
$ws3.Range("V17").Value = 1
$ws3.Range("A18").Value = "São Domingos"
$ws3.Range("B18").Value = 1
$ws3.Range("C18").Value = 0
$ws3.Range("D18").Value = 0
$ws3.Range("E18").Value = 1
$ws3.Range("F18").Value = 0
$ws3.Range("G18").Value = 1
$ws3.Range("H18").Value = 0
$ws3.Range("I18").Value = 0
$ws3.Range("J18").Value = 0
$ws3.Range("K18").Value = 0
$ws3.Range("L18").Value = 0
$ws3.Range("M18").Value = 0
$ws3.Range("N18").Value = 0
$ws3.Range("O18").Value = 0
$ws3.Range("P18").Value = 0
$ws3.Range("Q18").Value = 0
$ws3.Range("R18").Value = 0
$ws3.Range("S18").Value = 0
$ws3.Range("T18").Value = 0
$ws3.Range("U18").Value = 1
$ws3.Range("V18").Value = 0
$ws3.Range("A19").Value = "Carira"
$ws3.Range("B19").Value = 0
$ws3.Range("C19").Value = 1
$ws3.Range("D19").Value = 0
$ws3.Range("E19").Value = 0
$ws3.Range("F19").Value = 0
$ws3.Range("G19").Value = 1
$ws3.Range("H19").Value = 0
$ws3.Range("I19").Value = 1
$ws3.Range("J19").Value = 1
$ws3.Range("K19").Value = 0
$ws3.Range("L19").Value = 0
$ws3.Range("M19").Value = 1
$ws3.Range("N19").Value = 1
$ws3.Range("O19").Value = 0
$ws3.Range("P19").Value = 1
$ws3.Range("Q19").Value = 0
$ws3.Range("R19").Value = 0
$ws3.Range("S19").Value = 0
$ws3.Range("T19").Value = 0
$ws3.Range("U19").Value = 0
$ws3.Range("V19").Value = 0
$ws3.Range("A20").Value = "Areia Branca"
$ws3.Range("B20").Value = 1
$ws3.Range("C20").Value = 0
$ws3.Range("D20").Value = 0
$ws3.Range("E20").Value = 0
$ws3.Range("F20").Value = 1
$ws3.Range("G20").Value = 0
$ws3.Range("H20").Value = 1
$ws3.Range("I20").Value = 0
$ws3.Range("J20").Value = 0
$ws3.Range("K20").Value = 0
$ws3.Range("L20").Value = 0
$ws3.Range("M20").Value = 0
$ws3.Range("N20").Value = 0
$ws3.Range("O20").Value = 1
$ws3.Range("P20").Value = 0
$ws3.Range("Q20").Value = 1
$ws3.Range("R20").Value = 0
$ws3.Range("S20").Value = 0
$ws3.Range("T20").Value = 0
$ws3.Range("U20").Value = 0
$ws3.Range("V20").Value = 0
$ws3.Range("A21").Value = "Simão Dias"
$ws3.Range("B21").Value = 0
$ws3.Range("C21").Value = 0
$ws3.Range("D21").Value = 0
$ws3.Range("E21").Value = 1
$ws3.Range("F21").Value = 0
$ws3.Range("G21").Value = 0
$ws3.Range("H21").Value = 0
$ws3.Range("I21").Value = 0
$ws3.Range("J21").Value = 0
$ws3.Range("K21").Value = 0
$ws3.Range("L21").Value = 0
$ws3.Range("M21").Value = 1
$ws3.Range("N21").Value = 1
$ws3.Range("O21").Value = 0
$ws3.Range("P21").Value = 0
$ws3.Range("Q21").Value = 0
$ws3.Range("R21").Value = 0
$ws3.Range("S21").Value = 0
$ws3.Range("T21").Value = 0
$ws3.Range("U21").Value = 0
$ws3.Range("V21").Value = 0
$ws3.Range("A22").Value = "Aracaju"
$ws3.Range("B22").Value = 0
$ws3.Range("C22").Value = 0
$ws3.Range("D22").Value = 0
$ws3.Range("E22").Value = 0
$ws3.Range("F22").Value = 1
$ws3.Range("G22").Value = 0
$ws3.Range("H22").Value = 0
$ws3.Range("I22").Value = 0
$ws3.Range("J22").Value = 0
$ws3.Range("K22").Value = 0
$ws3.Range("L22").Value = 1
$ws3.Range("M22").Value = 0
$ws3.Range("N22").Value = 0
$ws3.Range("O22").Value = 1
$ws3.Range("P22").Value = 0
$ws3.Range("Q22").Value = 1
$ws3.Range("R22").Value = 0
$ws3.Range("S22").Value = 0
$ws3.Range("T22").Value = 1
$ws3.Range("U22").Value = 0
$ws3.Range("V22").Value = 0

# --- 7. Row heights that differ from the 15pt default (artifact of the paste) ---
$ws3.Rows.Item(1).RowHeight = 12.75
$ws3.Rows.Item(2).RowHeight = 14.25
$ws3.Rows.Item(4).RowHeight = 12.75
$ws3.Rows.Item(9).RowHeight = 16.5
$ws3.Rows.Item(10).RowHeight = 19.5
$ws3.Rows.Item(11).RowHeight = 15

# --- 8. Stray formatted-but-empty cells left over in column A below the table ---
$ws3.Range("A36").Value = ""
$ws3.Range("A37").Value = ""
$ws3.Range("A42").Value = ""
$ws3.Range("A43").Value = ""
$ws3.Range("A44").Value = ""
$ws3.Range("A45").Value = ""
$ws3.Range("A48").Value = ""
$ws3.Range("A49").Value = ""
$ws3.Range("A50").Value = ""
$ws3.Range("A51").Value = ""
$ws3.Range("A52").Value = ""
$ws3.Range("A53").Value = ""
$ws3.Range("A54").Value = ""
$ws3.Range("A55").Value = ""

# --- 9. Apply the shared label style (bold italic 10pt Georgia, dark green,
#        justify/center/wrap) to every labelled cell: header row, row labels,
#        and the stray empty cells below the table. ---
$labelCells = @("B1:V1","A2:A22") + ("A36","A37","A42","A43","A44","A45","A48","A49","A50","A51","A52","A53","A54","A55")
foreach ($addr in $labelCells) {
    $rng = $ws3.Range($addr)
    $rng.Font.Name = "Georgia"
    $rng.Font.FontStyle = "Bold Italic"
    $rng.Font.Size = 10
    $rng.Font.Color = 1463574
    $rng.HorizontalAlignment = -4130
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
}

# --- 10. Final selection/scroll position + make Plan3 the active/visible tab ---
$ws3.Application.Goto($ws3.Range("A1"), $false)
$ws3.Range("W8").Select()
$ws3.Activate()
